# "JoH paper under correction." — scores.xlsx update
#
# 1) Enable iterative-calculation bookkeeping (Formulas > Calculation Options):
#    the workbook gains a configured max-change of 1E-4 for iterative calc.
$excel.MaxChange = 0.0001

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2) Row 27 — average of the first score block (rows 16:25), mirroring the
#    B-alone / C:E-shared formula split already used elsewhere on the sheet.
$ws.Range("B27").Formula = "=AVERAGE(B16:B25)"
$ws.Range("C27:E27").Formula = "=AVERAGE(C16:C25)"

# 3) Row 53 — average of the second score block (rows 42:51).
$ws.Range("B53").Formula = "=AVERAGE(B42:B51)"
$ws.Range("C53:E53").Formula = "=AVERAGE(C42:C51)"

# 4) Row 79 — average of the third score block (rows 68:77).
$ws.Range("B79").Formula = "=AVERAGE(B68:B77)"
$ws.Range("C79:E79").Formula = "=AVERAGE(C68:C77)"

# 5) Move the saved selection/scroll position down to the newly added
#    row-79 summary (matches the author's last on-screen selection).
$ws.Range("B79:E79").Select() | Out-Null
